# Adds a new "Level" column (E) with header + per-row level numbers, and
# fixes a handful of casing/wording inconsistencies in column C (Level 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New column E: header "Level" (styled like the other headers in row 1)
#    plus an integer "Level" value for every data row (2-78).
# ---------------------------------------------------------------------

$ws.Range("E1").Value = "Level"
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").Borders.LineStyle = 1
$ws.Range("E1").HorizontalAlignment = -4108
$ws.Range("E1").VerticalAlignment = -4160

$levels = New-Object 'object[,]' 77,1
$levelValues = 1,1,1,1,1,1,1,1,1,1,1,1,1,2,2,2,1,1,1,1,2,1,1,1,2,2,2,1,1,1,1,1,1,1,1,1,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,3,3,2,2,2,2,1,2,2,2,2,2,2,2,2,2,2,2,2,1,2,2,2
for ($i = 0; $i -lt $levelValues.Length; $i++) {
    $levels[$i, 0] = $levelValues[$i]
}
$ws.Range("E2:E78").Value = $levels

# ---------------------------------------------------------------------
# 2) Column C (Level 2) text corrections - casing / wording tweaks.
# ---------------------------------------------------------------------

$cFixes = @{
    20 = "ICM Technology Handled"
    21 = "NonIVR Technology"
    22 = "Transferred calls"
    29 = "BCSC (Business Customer Service Center)"
    40 = "Billing > ViewBillCurrentPDF"
    41 = "Billing > ViewBillCurrentPDF historical"
    42 = "Usage and Rates > View Usage"
    43 = "Usage and Rates > Compare My Bills"
    44 = "Usage and Rates > Rate Comparison"
    45 = "Usage and Rates > Home Energy Checkup"
    46 = "Usage and Rates > Online Rate Enrollment"
    49 = "Outage > Subscribe Outage (EW Pages)"
    53 = "Payment Account NEW > CREATE_PAYMENT_ACCOUNT"
    54 = "Payment Account NEW > UPDATE_PAYMENT_ACCOUNT"
    55 = "Payment Account NEW > DELETE_PAYMENT_ACCOUNT"
    62 = "CARE/FERA"
    66 = "Alerts & Notifications > Go Paperless Alerts NEW"
    67 = "Alerts & Notifications > Notices & Services Information NEW"
    68 = "Alerts & Notifications > Event Day Alerts NEW"
    69 = "Alerts & Notifications > Service Visit Alerts NEW"
    70 = "Alerts & Notifications > Additional Communication Alerts NEW"
    71 = "Alerts & Notifications > Pay Plan & Shutoff Nonpayment Alerts NEW"
    72 = "Pilot Light Appointments  > SCHEDULE_SERVICE_APPOINTMENT"
    73 = "Pilot Light Appointments  > RESCHEDULE_SERVICE_APPOINTMENT"
    74 = "Pilot Light Appointments  > CANCEL_SERVICE_APPOINTMENT"
}

foreach ($row in $cFixes.Keys) {
    $ws.Cells.Item($row, 3).Value = $cFixes[$row]
}
